$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: rollNo f9 -> f8 (cascades into email + rollNo shared strings) ---
$ws.Range("C2").Value = "20l31a02f8@vignaniit.edu.in"
$ws.Range("D2").Value = "20l31a02f8"

# --- Row 3: drop this candidate (filtered out by cutOff) ---
# A3, D3, E3 had no hyperlink - fully clear them.
$ws.Range("A3").ClearContents()
$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()
# B3, C3, F3 carried hyperlinks + the "Hyperlink" style; clear the text but
# keep the styled (now-empty) cell, matching the target sheet shape.
$ws.Range("B3").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("F3").ClearContents()

# --- Hyperlinks: the runtime's Hyperlinks.Item(n).Delete() is a no-op, and
# Range(...).Hyperlinks.Delete() clears the *entire* sheet collection rather
# than just the addressed range, so rebuild the three we keep (B2/C2/F2)
# after wiping the lot, which drops the three dangling row-3 links (B3/C3/F3).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:aha@gmail.com", "", "", "aha@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:20l31a02f9@vignaniit.edu.in", "", "", "")
$ws.Hyperlinks.Add($ws.Range("F2"), "mailto:srithi@gmail.com", "", "", "")
# Re-adding a hyperlink stamps a fresh "Hyperlink" style slot; pin the cells
# back onto the workbook's existing Hyperlink cell style.
$ws.Range("B2").Style = "Hyperlink"
$ws.Range("C2").Style = "Hyperlink"
$ws.Range("F2").Style = "Hyperlink"

# --- View: was scrolled/selected around the (now gone) row-3 entry ---
[void]$ws.Range("A2").Select()
